# Updated symbol list on Sat Dec 24 13:55:40 UTC 2022 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values for the
# crypto ranking sheet. Column D holds numeric-looking quotes that must stay
# plain text (the sheet stores them as inline strings), so each assignment is
# prefixed with a literal apostrophe - exactly like typing '244.36 into an
# Excel cell - which keeps Excel from silently re-typing the cell as a
# number. Column E values are ordinary text and need no such treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D: Price -------------------------------------------------------
$ws.Range("D2").Value  = "'244.36"
$ws.Range("D3").Value  = "'21.91"
$ws.Range("D4").Value  = "'5.390"
$ws.Range("D5").Value  = "'0.05988"
$ws.Range("D7").Value  = "'0.8141"
$ws.Range("D8").Value  = "'0.9539"
$ws.Range("D9").Value  = "'0.1427"
$ws.Range("D10").Value = "'0.07432"
$ws.Range("D11").Value = "'0.03349"
$ws.Range("D12").Value = "'0.03053"
$ws.Range("D13").Value = "'0.09407"
$ws.Range("D14").Value = "'4.006"
$ws.Range("D15").Value = "'0.001589"
$ws.Range("D16").Value = "'0.04796"
$ws.Range("D17").Value = "'0.0005903"
$ws.Range("D18").Value = "'0.005109"
$ws.Range("D19").Value = "'0.005002"
$ws.Range("D20").Value = "'0.0009891"
$ws.Range("D22").Value = "'3.685"
$ws.Range("D23").Value = "'6.436"
$ws.Range("D26").Value = "'0.1332"
$ws.Range("D40").Value = "'0.03989"
$ws.Range("D41").Value = "'0.006563"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D44").Value = "'0.005764"
$ws.Range("D45").Value = "'0.00005277"
$ws.Range("D47").Value = "'0.9304"
$ws.Range("D48").Value = "'0.01239"

# --- Column E: Volume(1h) ---------------------------------------------------
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
